# Refresh crypto price/volume snapshot (GitHub Actions symbol-list update).
# Prices (col D) and 1h volume deltas (col E) are stored as plain text in
# this sheet, so each new value is entered with a leading apostrophe
# (quote-prefix) to force Excel to keep it as text instead of silently
# re-interpreting the numeric-looking string (e.g. "308.01") as a Number
# or the percentage strings (e.g. "-1.34%") as a Percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.01"
$ws.Range("E2").Value = "'-1.34%"

$ws.Range("D3").Value = "'36.05"
$ws.Range("E3").Value = "'-4.48%"

$ws.Range("D4").Value = "'5.119"
$ws.Range("E4").Value = "'-0.27%"

$ws.Range("D5").Value = "'0.07690"
$ws.Range("E5").Value = "'-2.73%"

$ws.Range("D6").Value = "'4.390"
$ws.Range("E6").Value = "'-0.44%"

$ws.Range("D7").Value = "'8.291"
$ws.Range("E7").Value = "'0.29%"

$ws.Range("D8").Value = "'1.842"
$ws.Range("E8").Value = "'-3.17%"

$ws.Range("E9").Value = "'-4.75%"

$ws.Range("D10").Value = "'0.9199"
$ws.Range("E10").Value = "'-0.23%"

$ws.Range("D11").Value = "'0.1110"
$ws.Range("E11").Value = "'-7.97%"

$ws.Range("D12").Value = "'0.1850"
$ws.Range("E12").Value = "'-4.44%"

$ws.Range("D13").Value = "'0.08718"
$ws.Range("E13").Value = "'-5.51%"

$ws.Range("D14").Value = "'0.03345"
$ws.Range("E14").Value = "'0.20%"

$ws.Range("D15").Value = "'0.09516"
$ws.Range("E15").Value = "'-1.07%"

$ws.Range("D16").Value = "'0.001384"
$ws.Range("E16").Value = "'-0.16%"

$ws.Range("D17").Value = "'0.006164"
$ws.Range("E17").Value = "'5.60%"

$ws.Range("E18").Value = "'-4.42%"

$ws.Range("E19").Value = "'-0.19%"

$ws.Range("E20").Value = "'19.30%"

$ws.Range("D21").Value = "'0.1291"
$ws.Range("E21").Value = "'1.49%"

$ws.Range("D23").Value = "'0.04338"
$ws.Range("E23").Value = "'-0.54%"

$ws.Range("E24").Value = "'-3.46%"

$ws.Range("D25").Value = "'0.004247"
$ws.Range("E25").Value = "'-1.49%"

$ws.Range("E26").Value = "'9.21%"

$ws.Range("D27").Value = "'0.0002904"

$ws.Range("D39").Value = "'0.02082"
$ws.Range("E39").Value = "'-2.21%"

$ws.Range("D40").Value = "'0.04907"
$ws.Range("E40").Value = "'-5.14%"

$ws.Range("D41").Value = "'0.007534"
$ws.Range("E41").Value = "'-1.43%"

$ws.Range("D42").Value = "'0.1346"
$ws.Range("E42").Value = "'-1.24%"

$ws.Range("D43").Value = "'0.008571"
$ws.Range("E43").Value = "'-5.88%"

$ws.Range("E44").Value = "'3.17%"

$ws.Range("D45").Value = "'0.008393"
$ws.Range("E45").Value = "'-2.39%"

$ws.Range("D46").Value = "'0.00006337"
$ws.Range("E46").Value = "'-5.40%"

$ws.Range("E47").Value = "'0.18%"

$ws.Range("D48").Value = "'0.003300"
$ws.Range("E48").Value = "'14.65%"

$ws.Range("E49").Value = "'20.45%"

$ws.Range("E50").Value = "'0.18%"

$ws.Range("E51").Value = "'0.18%"

Write-Output "Updated symbol list: 38 rows refreshed"
